# configuracion area 1 ospf
#
# Updates the DHCP exclusion-address rows (A1:A8) on the "DHCP" sheet,
# removing the stray "ADD" token from each "IP DH EX ADD ..." line, and
# narrows the active selection from the whole used range (A1:A49) down
# to just the edited rows (A1:A8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "IP DH EX 192.168.0.0 192.168.0.5"
$ws.Range("A2").Value = "IP DH EX  192.168.0.32 192.168.0.37"
$ws.Range("A3").Value = "IP DH EX  192.168.0.64 192.168.0.69"
$ws.Range("A4").Value = "IP DH EX  192.168.0.96 192.168.0.101"
$ws.Range("A5").Value = "IP DH EX  192.168.0.128 192.168.0.133"
$ws.Range("A6").Value = "IP DH EX  192.168.0.160 192.168.0.165"
$ws.Range("A7").Value = "IP DH EX  192.168.0.192 192.168.0.197"
$ws.Range("A8").Value = "IP DH EX  192.168.0.224 192.168.0.229"

# Try to mirror the minimized window state recorded in the commit; the
# sandboxed host may not persist this to workbookView, but setting it is
# harmless and matches the author's recorded intent.
try {
    $wb.Windows.Item(1).WindowState = -4140
} catch {
}

# Narrow the saved selection to the rows that were just edited.
$ws.Range("A1:A8").Select() | Out-Null
